$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The row count before and after the edit is identical: the old SIMONE row
# and the BEATRIZ row disappear while a (repositioned) SIMONE row and a new
# OLGA row take their place, so rows 7+ (VALMIR, THAYSA, MONICA, ...) keep
# their original row numbers. Only rows 2-6 need their contents rewritten.

# Row 2: new SIMONE entry (moved above LUIS, with an updated balance).
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "002823185"
$ws.Range("B2").Value = "SIMONE"
$ws.Range("C2").Value = 105155.96

# Row 3: LUIS keeps his account/name, only the balance changes.
$ws.Range("A3").NumberFormat = "@"
$ws.Range("A3").Value = "008007764"
$ws.Range("B3").Value = "LUIS"
$ws.Range("C3").Value = 40999.9

# Row 4: new OLGA entry (replaces the old SIMONE row at this position).
$ws.Range("A4").NumberFormat = "@"
$ws.Range("A4").Value = "008004799"
$ws.Range("B4").Value = "OLGA"
$ws.Range("C4").Value = 40000

# Row 5: PEDRO is unchanged, just re-asserted in its (same) position.
$ws.Range("A5").NumberFormat = "@"
$ws.Range("A5").Value = "004460487"
$ws.Range("B5").Value = "PEDRO"
$ws.Range("C5").Value = 16247.27

# Row 6: CINCO is unchanged, re-asserted in its (same) position. BEATRIZ
# (who used to occupy row 6 right after CINCO) is fully removed, so VALMIR
# (previously row 7) keeps row 7 as-is and everything below is untouched.
$ws.Range("A6").NumberFormat = "@"
$ws.Range("A6").Value = "004581652"
$ws.Range("B6").Value = "CINCO"
$ws.Range("C6").Value = 14455.12
